$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.305.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.24%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.321.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.12%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.62"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.46%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.13"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.70%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.43%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.318.93"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.23%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.03%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.50"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.04%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.05%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.23%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.68"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.55%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.272.67"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.30%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.734.62"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.12%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.314.16"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.46%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.58"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.08%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.13"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.75%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.65"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.72%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.64"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.97%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.29%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.79"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.77%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.172"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.04%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.03%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.18%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.59%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.53"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.50%  "

# Row 30
$ws.Range("B30").Value = "SuiNetwork"
$ws.Range("C30").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.35%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.73"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.31%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0729"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.00%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.93"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.35%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.83%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.97%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.32%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.04"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.26%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "322.43"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +9.80%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.91"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.11%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.40%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "138.18"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.02%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.83%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0940"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.39%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.20"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.76%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.30%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.560"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.66%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.48%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0213"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +17.12%  "

# Row 51
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.03"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.67%  "

